# Applies the weekly reshuffle of the Hortaliza / Vega Modelo de Temuco - Caigua
# dataset: rows 2-26 get their Fecha (D), Volumen (J), Precio minimo/maximo/promedio
# (K/L/M) and Precio $/Kg (P) columns permuted to new rows, per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as a group (everything else in the row is identical
# across all records, so only these columns actually change value).
$cols = @(4, 10, 11, 12, 13, 16)   # D, J, K, L, M, P

# Snapshot current values for rows 2..26 before any writes, since this is a
# permutation (source rows are also destination rows).
$snapshot = @{}
for ($r = 2; $r -le 26; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Mapping: destination row -> source row (values to copy into destination).
$mapping = @{
    2  = 3
    3  = 14
    4  = 18
    5  = 4
    6  = 17
    7  = 23
    8  = 6
    9  = 21
    10 = 8
    11 = 15
    12 = 16
    13 = 24
    14 = 5
    15 = 19
    16 = 13
    17 = 20
    18 = 7
    19 = 11
    20 = 26
    21 = 10
    22 = 22
    23 = 2
    24 = 25
    25 = 12
    26 = 9
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    if ($srcRow -eq $destRow) {
        continue
    }
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value2 = $srcVals[$c]
    }
}
